$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new profile row (row 90): Sanjib Panda / sanjib_panda, flag 0.
$ws.Range("A90").Value = "Sanjib Panda"
$ws.Range("B90").Value = "sanjib_panda"
$ws.Range("C90").Value = 0.0

# Match the existing formatting pattern used by the surrounding rows.
$ws.Range("A88").Copy()
$ws.Range("A90").PasteSpecial(-4122)

$ws.Range("B89").Copy()
$ws.Range("B90").PasteSpecial(-4122)

$ws.Range("C89").Copy()
$ws.Range("C90").PasteSpecial(-4122)
